$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark the next batch of MAUI tasks as complete -----------------------
# Rows 44-49 ("Create and implement logout page" .. "Add register page to
# shell nav") get their "Is Complete" cell (column D) set to "YES".
#
# The sheet has a live AutoFilter on B3:D62 that only shows blank "Is
# Complete" rows, so once these six rows are marked complete they drop out
# of the filtered view (rows 44-46 newly hidden; rows 47-49 were already
# hidden). The next chunk of previously-hidden rows (50-62) becomes the
# newly visible "next up" tasks.
#
# Temporarily unhide the whole 44:62 block before writing the values so
# the edit doesn't leave stray auto-row-height markers on rows that are
# hidden at edit time, then re-apply the correct hidden state afterwards.
$ws.Range("44:62").EntireRow.Hidden = $false

$ws.Range("D44:D49").Value = "YES"

$ws.Range("44:49").EntireRow.Hidden = $true
$ws.Range("50:62").EntireRow.Hidden = $false

# --- Clear the leftover fill on the rows that just got checked off -------
# B44:C46 previously carried the "upcoming" light-blue highlight and
# B47:C49 the "upcoming" yellow highlight; once the task is complete that
# highlight is cleared back to no fill (matching the rest of the
# completed/checked-off rows above).
$ws.Range("B44:C46").Interior.Pattern = -4142
$ws.Range("B44:C46").Copy() | Out-Null
$ws.Range("B47:C49").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# D44:D46 swap from the "upcoming" highlighted YES style to the plain
# centered YES style used elsewhere (no fill).
$ws.Range("D44:D46").Interior.Pattern = -4142

# --- Restore cursor/scroll position ---------------------------------------
$ws.Range("B51").Select()
